$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'322.15"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-2.81%"
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'42.80"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-5.80%"
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").Value = "'5.146"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-8.27%"
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'0.08187"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-1.98%"
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'4.277"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-3.63%"
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").Value = "'1.800"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-13.79%"
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'0.9316"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-3.51%"
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").Value = "'0.1110"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-5.26%"
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("E10").Value = "'-2.99%"
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").Value = "'0.09455"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-4.07%"
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'0.04635"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'0.38%"
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").Value = "'7.397"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-28.63%"
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("E14").Value = "'-0.20%"
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'0.001293"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-0.27%"
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'0.005762"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-5.07%"
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'3.364"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-0.39%"
$ws.Range("E17").ClearFormats()

# Row 19
$ws.Range("D19").Value = "'0.3374"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'0.95%"
$ws.Range("E19").ClearFormats()

# Row 21
$ws.Range("E21").Value = "'-12.44%"
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").Value = "'0.04160"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-0.58%"
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'0.001246"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-5.30%"
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'0.004453"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-2.23%"
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("E25").Value = "'-7.85%"
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").Value = "'-20.52%"
$ws.Range("E26").ClearFormats()

# Row 38
$ws.Range("E38").Value = "'1.10%"
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'0.05587"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-2.97%"
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").Value = "'0.007923"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'0.71%"
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").Value = "'0.1398"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-2.52%"
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'0.006545"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-9.96%"
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").Value = "'0.002041"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'0.82%"
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = "'0.008346"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-8.20%"
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'0.3496"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-1.44%"
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").Value = "'0.00006969"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-2.26%"
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").Value = "'0.003473"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'-1.50%"
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("D49").Value = "'0.003532"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'0.69%"
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").ClearFormats()
